$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.559.43"
$ws.Range("E2").Value = "'  +0.24%  "
$ws.Range("D3").Value = "'1.846.47"
$ws.Range("E3").Value = "'  -0.20%  "
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'264.24"
$ws.Range("E5").Value = "'  +1.50%  "
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  -0.01%  "
$ws.Range("D7").Value = "'0.5212"
$ws.Range("E7").Value = "'  +1.08%  "
$ws.Range("D8").Value = "'0.3237"
$ws.Range("E8").Value = "'  -0.43%  "
$ws.Range("D9").Value = "'0.06807"
$ws.Range("E9").Value = "'  +0.72%  "
$ws.Range("D10").Value = "'18.78"
$ws.Range("E10").Value = "'  -0.54%  "
$ws.Range("D11").Value = "'0.7782"
$ws.Range("E11").Value = "'  +0.85%  "
$ws.Range("D12").Value = "'0.07771"
$ws.Range("E12").Value = "'  +0.74%  "
$ws.Range("D13").Value = "'1.864.24"
$ws.Range("E13").Value = "'  +0.79%  "
$ws.Range("D14").Value = "'88.51"
$ws.Range("E14").Value = "'  -0.14%  "
$ws.Range("D15").Value = "'5.030"
$ws.Range("E15").Value = "'  +0.03%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "'  -0.07%  "
$ws.Range("D17").Value = "'14.00"
$ws.Range("E17").Value = "'  -0.65%  "
$ws.Range("E18").Value = "'  -0.01%  "
$ws.Range("D19").Value = "'0.000007962"
$ws.Range("E19").Value = "'  +0.67%  "
$ws.Range("D20").Value = "'26.601.43"
$ws.Range("E20").Value = "'  +0.36%  "
$ws.Range("D21").Value = "'4.632"
$ws.Range("E21").Value = "'  +2.37%  "
$ws.Range("D22").Value = "'9.462"
$ws.Range("E22").Value = "'  -0.78%  "
$ws.Range("D23").Value = "'6.016"
$ws.Range("E23").Value = "'  +1.60%  "
$ws.Range("D24").Value = "'143.42"
$ws.Range("E24").Value = "'  -0.72%  "
$ws.Range("D25").Value = "'2.176"
$ws.Range("E25").Value = "'  -7.48%  "
$ws.Range("D26").Value = "'1.678"
$ws.Range("E26").Value = "'  +1.73%  "
$ws.Range("D27").Value = "'17.01"
$ws.Range("E27").Value = "'  +0.17%  "
$ws.Range("D28").Value = "'111.76"
$ws.Range("E28").Value = "'  +0.55%  "
$ws.Range("D29").Value = "'4.192"
$ws.Range("E29").Value = "'  -0.35%  "
$ws.Range("D30").Value = "'0.08747"
$ws.Range("E30").Value = "'  -0.12%  "
$ws.Range("D31").Value = "'4.116"
$ws.Range("E31").Value = "'  -1.36%  "
$ws.Range("D32").Value = "'0.04842"
$ws.Range("E32").Value = "'  +0.57%  "
$ws.Range("D33").Value = "'0.7215"
$ws.Range("E33").Value = "'  +4.91%  "
$ws.Range("D34").Value = "'1.131"
$ws.Range("E34").Value = "'  -0.33%  "
$ws.Range("D35").Value = "'2.860"
$ws.Range("E35").Value = "'  +0.69%  "
$ws.Range("D36").Value = "'3.105"
$ws.Range("E36").Value = "'  -0.22%  "
$ws.Range("D37").Value = "'0.01793"
$ws.Range("E37").Value = "'  -0.49%  "
$ws.Range("D38").Value = "'2.219"
$ws.Range("E38").Value = "'  +0.36%  "
$ws.Range("D39").Value = "'0.4862"
$ws.Range("E39").Value = "'  -0.79%  "
$ws.Range("D40").Value = "'111.13"
$ws.Range("E40").Value = "'  -1.72%  "
$ws.Range("D41").Value = "'0.8943"
$ws.Range("E41").Value = "'  -0.53%  "
$ws.Range("D42").Value = "'6.033"
$ws.Range("E42").Value = "'  -1.63%  "
$ws.Range("D43").Value = "'1.0000"
$ws.Range("E43").Value = "'  -0.03%  "
$ws.Range("D44").Value = "'7.625"
$ws.Range("E44").Value = "'  -1.93%  "
$ws.Range("D45").Value = "'0.4213"
$ws.Range("E45").Value = "'  -0.17%  "
$ws.Range("D46").Value = "'0.05891"
$ws.Range("E46").Value = "'  -0.01%  "
$ws.Range("D47").Value = "'9.070"
$ws.Range("E47").Value = "'  +0.02%  "
$ws.Range("D48").Value = "'0.1239"
$ws.Range("E48").Value = "'  -1.71%  "
$ws.Range("D49").Value = "'35.05"
$ws.Range("E49").Value = "'  -0.20%  "
$ws.Range("E50").Value = "'  +3.73%  "
$ws.Range("D51").Value = "'59.95"
$ws.Range("E51").Value = "'  +1.38%  "
